$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.8712260127067566
$ws.Range("C2").Value2 = 0.9262280464172363
$ws.Range("D2").Value2 = 0.9063977599143982
$ws.Range("E2").Value2 = 0.8978849649429321
$ws.Range("F2").Value2 = 0.9083855152130127
$ws.Range("B3").Value2 = 0.6321839094161987
$ws.Range("C3").Value2 = 0.2801358103752136
$ws.Range("D3").Value2 = 0.9590293169021606
$ws.Range("E3").Value2 = 0.3882348537445068
$ws.Range("F3").Value2 = 0.6361019611358643
$ws.Range("B4").Value2 = 0.1084337383508682
$ws.Range("C4").Value2 = 0.0476190485060215
$ws.Range("D4").Value2 = 0.9799873828887939
$ws.Range("E4").Value2 = 0.06617604941129684
$ws.Range("F4").Value2 = 0.5208502411842346
$ws.Range("B5").Value2 = 0.8858789801597595
$ws.Range("C5").Value2 = 0.8713151812553406
$ws.Range("D5").Value2 = 0.8660573363304138
$ws.Range("E5").Value2 = 0.8785362243652344
$ws.Range("F5").Value2 = 0.8653949499130249
$ws.Range("B6").Value2 = 0.7549019455909729
$ws.Range("C6").Value2 = 0.8271835446357727
$ws.Range("D6").Value2 = 0.9255436658859253
$ws.Range("E6").Value2 = 0.7893911004066467
$ws.Range("F6").Value2 = 0.8863431811332703
$ws.Range("B7").Value2 = 0.6270411610603333
$ws.Range("C7").Value2 = 0.7082257270812988
$ws.Range("D7").Value2 = 0.8476993441581726
$ws.Range("E7").Value2 = 0.6651649475097656
$ws.Range("F7").Value2 = 0.7969042062759399
$ws.Range("B8").Value2 = 0.4575389921665192
$ws.Range("C8").Value2 = 0.2573099434375763
$ws.Range("D8").Value2 = 0.9153009653091431
$ws.Range("E8").Value2 = 0.3293819427490234
$ws.Range("F8").Value2 = 0.6152399182319641
$ws.Range("B9").Value2 = 0.9366075396537781
$ws.Range("C9").Value2 = 0.973239541053772
$ws.Range("D9").Value2 = 0.9149070382118225
$ws.Range("E9").Value2 = 0.9545717835426331
$ws.Range("F9").Value2 = 0.6148869395256042
$ws.Range("B10").Value2 = 0.8953509330749512
$ws.Range("C10").Value2 = 0.92912358045578
$ws.Range("D10").Value2 = 0.8481720685958862
$ws.Range("E10").Value2 = 0.9119241833686829
$ws.Range("F10").Value2 = 0.6663520932197571
$ws.Range("B11").Value2 = 0.3322683572769165
$ws.Range("C11").Value2 = 0.1522693932056427
$ws.Range("D11").Value2 = 0.9379136562347412
$ws.Range("E11").Value2 = 0.2088348865509033
$ws.Range("F11").Value2 = 0.5674328804016113
$ws.Range("B12").Value2 = 0.1022727265954018
$ws.Range("C12").Value2 = 0.02564102597534657
$ws.Range("D12").Value2 = 0.9668294787406921
$ws.Range("E12").Value2 = 0.04100195690989494
$ws.Range("F12").Value2 = 0.5096198320388794
$ws.Range("B13").Value2 = 0.4586597084999084
$ws.Range("C13").Value2 = 0.5759562849998474
$ws.Range("D13").Value2 = 0.9204223155975342
$ws.Range("E13").Value2 = 0.5106584429740906
$ws.Range("F13").Value2 = 0.7615707516670227
$ws.Range("B14").Value2 = 0.5447154641151428
$ws.Range("C14").Value2 = 0.07854630798101425
$ws.Range("D14").Value2 = 0.9336590170860291
$ws.Range("E14").Value2 = 0.1372948586940765
$ws.Range("F14").Value2 = 0.5369080901145935
$ws.Range("B15").Value2 = 0.8745366334915161
$ws.Range("C15").Value2 = 0.7177627086639404
$ws.Range("D15").Value2 = 0.8703120350837708
$ws.Range("E15").Value2 = 0.7884314656257629
$ws.Range("F15").Value2 = 0.8327499628067017
$ws.Range("B16").Value2 = 0.8141592741012573
$ws.Range("C16").Value2 = 0.7101754546165466
$ws.Range("D16").Value2 = 0.9492594003677368
$ws.Range("E16").Value2 = 0.7586202621459961
$ws.Range("F16").Value2 = 0.8448365926742554
$ws.Range("B17").Value2 = 0.973128616809845
$ws.Range("C17").Value2 = 0.8718830347061157
$ws.Range("D17").Value2 = 0.9860541820526123
$ws.Range("E17").Value2 = 0.9197273850440979
$ws.Range("F17").Value2 = 0.9347271919250488
$ws.Range("B18").Value2 = 0.787106454372406
$ws.Range("C18").Value2 = 0.6126021146774292
$ws.Range("D18").Value2 = 0.962653636932373
$ws.Range("E18").Value2 = 0.6889759302139282
$ws.Range("F18").Value2 = 0.8003019094467163
$ws.Range("B19").Value2 = 0.6955307126045227
$ws.Range("C19").Value2 = 0.6505551934242249
$ws.Range("D19").Value2 = 0.9234951138496399
$ws.Range("E19").Value2 = 0.672291100025177
$ws.Range("F19").Value2 = 0.8057453036308289
$ws.Range("B20").Value2 = 0.7321428656578064
$ws.Range("C20").Value2 = 0.4324894547462463
$ws.Range("D20").Value2 = 0.8915852308273315
$ws.Range("E20").Value2 = 0.5437661409378052
$ws.Range("F20").Value2 = 0.7023507356643677
$ws.Range("B21").Value2 = 0.8400900959968567
$ws.Range("C21").Value2 = 0.5088676810264587
$ws.Range("D21").Value2 = 0.966041624546051
$ws.Range("E21").Value2 = 0.6338143348693848
$ws.Range("F21").Value2 = 0.7514653205871582
$ws.Range("B22").Value2 = 0.7851351499557495
$ws.Range("C22").Value2 = 0.6755813956260681
$ws.Range("D22").Value2 = 0.9309801459312439
$ws.Range("E22").Value2 = 0.7262495160102844
$ws.Range("F22").Value2 = 0.8232992887496948
$ws.Range("B23").Value2 = 0.320277214050293
$ws.Range("C23").Value2 = 0.9582155346870422
$ws.Range("D23").Value2 = 0.4795146584510803
$ws.Range("E23").Value2 = 0.4800877869129181
$ws.Range("F23").Value2 = 0.6387460231781006
$ws.Range("B24").Value2 = 0
$ws.Range("C24").Value2 = 0
$ws.Range("D24").Value2 = 0.9886542558670044
$ws.Range("E24").Value2 = 0
$ws.Range("F24").Value2 = 0.5
$ws.Range("B25").Value2 = 0
$ws.Range("C25").Value2 = 0
$ws.Range("D25").Value2 = 0.7605578303337097
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 0.5
$ws.Range("B26").Value2 = 0.7305143475532532
$ws.Range("C26").Value2 = 0.5026053190231323
$ws.Range("D26").Value2 = 0.7522060871124268
$ws.Range("E26").Value2 = 0.5954979062080383
$ws.Range("F26").Value2 = 0.6984953284263611
$ws.Range("B27").Value2 = 0.9261554479598999
$ws.Range("C27").Value2 = 0.8891918659210205
$ws.Range("D27").Value2 = 0.9422470927238464
$ws.Range("E27").Value2 = 0.907296895980835
$ws.Range("F27").Value2 = 0.9280794262886047
$ws.Range("B28").Value2 = 0.891465425491333
$ws.Range("C28").Value2 = 0.9216910004615784
$ws.Range("D28").Value2 = 0.8941853046417236
$ws.Range("E28").Value2 = 0.9063258171081543
$ws.Range("F28").Value2 = 0.8907586336135864
$ws.Range("B29").Value2 = 0.7517961859703064
$ws.Range("C29").Value2 = 0.718028724193573
$ws.Range("D29").Value2 = 0.9344468712806702
$ws.Range("E29").Value2 = 0.7345240712165833
$ws.Range("F29").Value2 = 0.8418802618980408
$ws.Range("B30").Value2 = 0.5366747975349426
$ws.Range("C30").Value2 = 0.5028637051582336
$ws.Range("D30").Value2 = 0.9359439015388489
$ws.Range("E30").Value2 = 0.519218921661377
$ws.Range("F30").Value2 = 0.7353983521461487
$ws.Range("B31").Value2 = 0.3299180269241333
$ws.Range("C31").Value2 = 0.574999988079071
$ws.Range("D31").Value2 = 0.9648597836494446
$ws.Range("E31").Value2 = 0.4192703664302826
$ws.Range("F31").Value2 = 0.774327278137207
$ws.Range("B32").Value2 = 0.4123989343643188
$ws.Range("C32").Value2 = 0.3903061151504517
$ws.Range("D32").Value2 = 0.9639930725097656
$ws.Range("E32").Value2 = 0.4010479748249054
$ws.Range("F32").Value2 = 0.6862912774085999
$ws.Range("B33").Value2 = 0.5658436417579651
$ws.Range("C33").Value2 = 0.3450439274311066
$ws.Range("D33").Value2 = 0.9422470927238464
$ws.Range("E33").Value2 = 0.4286822974681854
$ws.Range("F33").Value2 = 0.6636527180671692
$ws.Range("B34").Value2 = 0.3682926893234253
$ws.Range("C34").Value2 = 0.4026666581630707
$ws.Range("D34").Value2 = 0.9619445204734802
$ws.Range("E34").Value2 = 0.3847128450870514
$ws.Range("F34").Value2 = 0.6908193826675415
$ws.Range("B35").Value2 = 0.7723866105079651
$ws.Range("C35").Value2 = 0.6938341856002808
$ws.Range("D35").Value2 = 0.8864639401435852
$ws.Range("E35").Value2 = 0.7310056686401367
$ws.Range("F35").Value2 = 0.8176870942115784
$ws.Range("B36").Value2 = 0.5500413775444031
$ws.Range("C36").Value2 = 0.5657167434692383
$ws.Range("D36").Value2 = 0.833832323551178
$ws.Range("E36").Value2 = 0.5577684640884399
$ws.Range("F36").Value2 = 0.7302522659301758
$ws.Range("B37").Value2 = 0
$ws.Range("C37").Value2 = 0
$ws.Range("D37").Value2 = 0.9913331270217896
$ws.Range("E37").Value2 = 0
$ws.Range("F37").Value2 = 0.5
$ws.Range("B38").Value2 = 0.8029423356056213
$ws.Range("C38").Value2 = 0.885189950466156
$ws.Range("D38").Value2 = 0.8774030804634094
$ws.Range("E38").Value2 = 0.8420620560646057
$ws.Range("F38").Value2 = 0.8790176510810852
$ws.Range("B39").Value2 = 0.82596355676651
$ws.Range("C39").Value2 = 0.8098770976066589
$ws.Range("D39").Value2 = 0.7548061609268188
$ws.Range("E39").Value2 = 0.8178407549858093
$ws.Range("F39").Value2 = 0.7239252328872681
$ws.Range("B40").Value2 = 0.5933352112770081
$ws.Range("C40").Value2 = 0.6163097620010376
$ws.Range("D40").Value2 = 0.7834856510162354
$ws.Range("E40").Value2 = 0.6046038269996643
$ws.Range("F40").Value2 = 0.730593740940094
$ws.Range("B41").Value2 = 0.5440115332603455
$ws.Range("C41").Value2 = 0.5791090726852417
$ws.Range("D41").Value2 = 0.953514039516449
$ws.Range("E41").Value2 = 0.5610114932060242
$ws.Range("F41").Value2 = 0.7764326930046082
$ws.Range("B42").Value2 = 0.588320791721344
$ws.Range("C42").Value2 = 0.5980295538902283
$ws.Range("D42").Value2 = 0.737551212310791
$ws.Range("E42").Value2 = 0.5931349992752075
$ws.Range("F42").Value2 = 0.7006019353866577
$ws.Range("B43").Value2 = 0.4675492346286774
$ws.Range("C43").Value2 = 0.6156462430953979
$ws.Range("D43").Value2 = 0.7988496422767639
$ws.Range("E43").Value2 = 0.5314732193946838
$ws.Range("F43").Value2 = 0.7280842065811157
$ws.Range("B44").Value2 = 0.5450361371040344
$ws.Range("C44").Value2 = 0.2902660965919495
$ws.Range("D44").Value2 = 0.7857705354690552
$ws.Range("E44").Value2 = 0.3787977695465088
$ws.Range("F44").Value2 = 0.6099561452865601
$ws.Range("B45").Value2 = 0.6148973703384399
$ws.Range("C45").Value2 = 0.7100643515586853
$ws.Range("D45").Value2 = 0.8021588325500488
$ws.Range("E45").Value2 = 0.6590626835823059
$ws.Range("F45").Value2 = 0.7730826139450073
$ws.Range("B46").Value2 = 0.3636363446712494
$ws.Range("C46").Value2 = 0.1176470592617989
$ws.Range("D46").Value2 = 0.9970847964286804
$ws.Range("E46").Value2 = 0.1777773946523666
$ws.Range("F46").Value2 = 0.5585470199584961
$ws.Range("B47").Value2 = 0.935258150100708
$ws.Range("C47").Value2 = 0.8556617498397827
$ws.Range("D47").Value2 = 0.8220926523208618
$ws.Range("E47").Value2 = 0.8936906456947327
$ws.Range("F47").Value2 = 0.7225183844566345
$ws.Range("B48").Value2 = 0.3549663722515106
$ws.Range("C48").Value2 = 0.6056718230247498
$ws.Range("D48").Value2 = 0.8255594372749329
$ws.Range("E48").Value2 = 0.4476043283939362
$ws.Range("F48").Value2 = 0.7301394939422607
$ws.Range("B49").Value2 = 0
$ws.Range("C49").Value2 = 0
$ws.Range("D49").Value2 = 0.9928301572799683
$ws.Range("E49").Value2 = 0
$ws.Range("F49").Value2 = 0.4997620284557343
